# Apply the "Add data for 2022-10-18" update to the carjacking workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-10-10"

# Update the October row label to reflect the new "through" date.
$ws.Range("A11").Value = "October (through 10-10)"

# Update the October (row 11) year-by-year counts.
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 17
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = 25
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = 43
$ws.Range("H11").Value = 67
$ws.Range("I11").Value = 35

# Update the Total (row 12) year-by-year counts.
$ws.Range("B12").Value = 234
$ws.Range("C12").Value = 446
$ws.Range("D12").Value = 644
$ws.Range("E12").Value = 573
$ws.Range("F12").Value = 433
$ws.Range("G12").Value = 944
$ws.Range("H12").Value = 1314
$ws.Range("I12").Value = 1313
